$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Receptor-expressing cells dropped from 2 to 1, recalculated downstream metrics)
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.036942
$ws.Range("N2").Value = 0.110826
$ws.Range("O2").Value = 0.02099032928903418
$ws.Range("P2").Value = 0.02099032928903418
$ws.Range("Q2").Value = 0.004948380899999999
$ws.Range("R2").Value = 0.0445354281
$ws.Range("S2").Value = 0.02099032928903418
$ws.Range("T2").Value = 0.02099032928903418

# Row 3 updates (specificity values recalculated due to change in row 2's receptor values)
$ws.Range("O3").Value = 0.5358731102718634
$ws.Range("P3").Value = 0.5358731102718634
$ws.Range("S3").Value = 0.5358731102718634
$ws.Range("T3").Value = 0.5358731102718634

# Row 4 updates (specificity values recalculated due to change in row 2's receptor values)
$ws.Range("O4").Value = 0.4431365604391025
$ws.Range("P4").Value = 0.4431365604391026
$ws.Range("S4").Value = 0.4431365604391025
$ws.Range("T4").Value = 0.4431365604391026
